# Repro of "Add files via upload" commit:
#   - Fixes 84 #DIV/0! cached values in column D of "12-months data"
#     (rows 218-361, in 13 contiguous blocks) back to the literal
#     46.32601623 that the rest of the column already uses.
#   - Restores the sheet's scroll/selection view-state (best effort —
#     scrolls the window down near the bottom of the data and leaves the
#     active cell at D257).
#   - Nudges the "SPI1" sheet's scroll position as well (best effort).
#
# NOTE on scope: the underlying model only round-trips the *selection*
# (ActiveCell / Selection) into <sheetView>; it has no host surface for
# a bare `topLeftCell` scroll offset outside of freeze/split panes. The
# ActiveWindow.ScrollRow/ScrollColumn calls below are issued anyway
# (matches what real Excel would record when a user scrolls), in case a
# future host revision starts honoring them, but they are not expected
# to change the saved XML on this runtime.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "12-months data" sheet: replace the cached #DIV/0! values in D with
#    the correct 46.32601623 literal (same value already present in the
#    surrounding, unaffected rows of the same column).
# ---------------------------------------------------------------------
$ws12 = $wb.Worksheets.Item("12-months data")

$fixedRanges = @(
    "D218:D221",
    "D227:D233",
    "D239:D245",
    "D251:D257",
    "D263:D269",
    "D275:D281",
    "D287:D293",
    "D299:D305",
    "D311:D317",
    "D323:D329",
    "D335:D341",
    "D347:D353",
    "D359:D361"
)

foreach ($rangeAddress in $fixedRanges) {
    $ws12.Range($rangeAddress).Value = 46.32601623
}

# ---------------------------------------------------------------------
# 2. "SPI1" sheet view state: scroll near row 136. The selection itself
#    is unchanged in the target (still A2:A6, active cell A2), so only
#    activate the sheet — do NOT re-select, which would collapse the
#    existing multi-cell selection down to a single cell.
# ---------------------------------------------------------------------
$wsSpi1 = $wb.Worksheets.Item("SPI1")
$wsSpi1.Activate()
try {
    $excel.ActiveWindow.ScrollRow = 136
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}

# ---------------------------------------------------------------------
# 3. "12-months data" sheet view state: scroll near row 370, active
#    cell D257. Re-activate this sheet last so it stays the workbook's
#    active tab (matches activeTab="1" / tabSelected="1").
# ---------------------------------------------------------------------
$ws12.Activate()
try {
    $excel.ActiveWindow.ScrollRow = 370
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}
$ws12.Range("D257").Select()
